# Update countries & provincias Spain
# Refresh COVID-19 country stats (totals/new/active/recovered/critical/deaths)
# and bump the "last updated" timestamp. A handful of countries swap rank
# (adjacent row pairs) because the table is sorted by total cases.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Header timestamp
$ws.Range("A1").Value = 'Datos actualizados a 24 de Julio de 2020 a las 18:53'

# Row 4: Estados Unidos
$ws.Range("B4").Value = 4200751
$ws.Range("C4").Value = 30433
$ws.Range("D4").Value = 1988921
$ws.Range("E4").Value = 2064142
$ws.Range("G4").Value = 339
$ws.Range("H4").Value = 147688

# Row 5: Brasil
$ws.Range("B5").Value = 2303661
$ws.Range("C5").Value = 13710
$ws.Range("E5").Value = 648984
$ws.Range("G5").Value = 233
$ws.Range("H5").Value = 84440

# Row 6: India
$ws.Range("B6").Value = 1333553
$ws.Range("C6").Value = 45423
$ws.Range("D6").Value = 847803
$ws.Range("E6").Value = 454360
$ws.Range("G6").Value = 745
$ws.Range("H6").Value = 31390

# Row 19: Turquia
$ws.Range("B19").Value = 224252
$ws.Range("C19").Value = 937
$ws.Range("D19").Value = 207374
$ws.Range("E19").Value = 11298
$ws.Range("G19").Value = 17
$ws.Range("H19").Value = 5580

# Row 21: Alemania
$ws.Range("B21").Value = 205402
$ws.Range("C21").Value = 260
$ws.Range("E21").Value = 6811
$ws.Range("G21").Value = 4
$ws.Range("H21").Value = 9191

# Row 40: Israel
$ws.Range("B40").Value = 59475
$ws.Range("C40").Value = 1493
$ws.Range("D40").Value = 26797
$ws.Range("E40").Value = 32230
$ws.Range("G40").Value = 6
$ws.Range("H40").Value = 448

# Row 41: Republica Dominicana
$ws.Range("A41").Value = 'Republica Dominicana'
$ws.Range("B41").Value = 59077
$ws.Range("C41").Value = 1462
$ws.Range("D41").Value = 27625
$ws.Range("E41").Value = 30416
$ws.Range("G41").Value = 30
$ws.Range("H41").Value = 1036

# Row 42: Emiratos Arabes Unidos
$ws.Range("A42").Value = 'Emiratos Arabes Unidos'
$ws.Range("B42").Value = 58249
$ws.Range("C42").Value = 261
$ws.Range("D42").Value = 51235
$ws.Range("E42").Value = 6671
$ws.Range("G42").Value = 1
$ws.Range("H42").Value = 343

# Row 45: Portugal
$ws.Range("B45").Value = 49692
$ws.Range("C45").Value = 313
$ws.Range("D45").Value = 34687
$ws.Range("E45").Value = 13293
$ws.Range("G45").Value = 7
$ws.Range("H45").Value = 1712

# Row 47: Guatemala
$ws.Range("A47").Value = 'Guatemala'
$ws.Range("B47").Value = 43283
$ws.Range("C47").Value = 1091
$ws.Range("D47").Value = 30150
$ws.Range("E47").Value = 11464
$ws.Range("G47").Value = 37
$ws.Range("H47").Value = 1669

# Row 48: Rumania
$ws.Range("A48").Value = 'Rumania'
$ws.Range("B48").Value = 42394
$ws.Range("C48").Value = 1119
$ws.Range("D48").Value = 25349
$ws.Range("E48").Value = 14895
$ws.Range("G48").Value = 24
$ws.Range("H48").Value = 2150

# Row 60: Argelia
$ws.Range("A60").Value = 'Argelia'
$ws.Range("B60").Value = 26159
$ws.Range("C60").Value = 675
$ws.Range("D60").Value = 17369
$ws.Range("E60").Value = 7654
$ws.Range("G60").Value = 12
$ws.Range("H60").Value = 1136

# Row 61: Irlanda
$ws.Range("A61").Value = 'Irlanda'
$ws.Range("B61").Value = 25826
$ws.Range("D61").Value = 23364
$ws.Range("E61").Value = 699
$ws.Range("H61").Value = 1763

# Row 65: Uzbekistan
$ws.Range("B65").Value = 19179
$ws.Range("C65").Value = 311
$ws.Range("D65").Value = 10203
$ws.Range("E65").Value = 8872
$ws.Range("G65").Value = 2
$ws.Range("H65").Value = 104

# Row 69: Kenia
$ws.Range("B69").Value = 16268
$ws.Range("C69").Value = 667
$ws.Range("D69").Value = 7446
$ws.Range("E69").Value = 8548
$ws.Range("G69").Value = 11
$ws.Range("H69").Value = 274

# Row 71: Chequia
$ws.Range("B71").Value = 14924
$ws.Range("C71").Value = 124
$ws.Range("D71").Value = 9402
$ws.Range("E71").Value = 5154
$ws.Range("G71").Value = 3
$ws.Range("H71").Value = 368

# Row 78: Etiopia
$ws.Range("B78").Value = 12693
$ws.Range("C78").Value = 760
$ws.Range("D78").Value = 5785
$ws.Range("E78").Value = 6708
$ws.Range("G78").Value = 3
$ws.Range("H78").Value = 200

# Row 90: Guayana Francesa
$ws.Range("A90").Value = 'Guayana Francesa'
$ws.Range("B90").Value = 7251
$ws.Range("C90").Value = 165
$ws.Range("D90").Value = 5522
$ws.Range("E90").Value = 1688
$ws.Range("G90").Value = 1
$ws.Range("H90").Value = 41

# Row 91: Haiti
$ws.Range("A91").Value = 'Haiti'
$ws.Range("B91").Value = 7197
$ws.Range("C91").Value = 30
$ws.Range("D91").Value = 4236
$ws.Range("E91").Value = 2807
$ws.Range("H91").Value = 154

# Row 92: Tayikistan
$ws.Range("A92").Value = 'Tayikistan'
$ws.Range("B92").Value = 7104
$ws.Range("C92").Value = 44
$ws.Range("D92").Value = 5851
$ws.Range("E92").Value = 1195
$ws.Range("H92").Value = 58

# Row 96: Luxemburgo
$ws.Range("B96").Value = 6056
$ws.Range("C96").Value = 104
$ws.Range("D96").Value = 4647
$ws.Range("E96").Value = 1297

# Row 97: Republica de Yibuti
$ws.Range("B97").Value = 5039
$ws.Range("C97").Value = 8
$ws.Range("D97").Value = 4949
$ws.Range("E97").Value = 32

# Row 114: Sri Lanka
$ws.Range("B114").Value = 2764
$ws.Range("C114").Value = 11
$ws.Range("E114").Value = 659

# Row 117: Cuba
$ws.Range("B117").Value = 2469
$ws.Range("C117").Value = 3
$ws.Range("D117").Value = 2341
$ws.Range("E117").Value = 41

# Row 135: Mozambique
$ws.Range("B135").Value = 1590
$ws.Range("C135").Value = 8
$ws.Range("D135").Value = 532
$ws.Range("E135").Value = 1047

# Row 137: Tunez
$ws.Range("B137").Value = 1425
$ws.Range("C137").Value = 19
$ws.Range("D137").Value = 1124
$ws.Range("E137").Value = 251

# Row 149: Principado de Andorra
$ws.Range("B149").Value = 897
$ws.Range("C149").Value = 8
$ws.Range("E149").Value = 42

# Row 158: Reunion
$ws.Range("B158").Value = 654
$ws.Range("C158").Value = 8
$ws.Range("E158").Value = 93

# Row 181: Trinidad yTobago
$ws.Range("B181").Value = 142
$ws.Range("C181").Value = 1
$ws.Range("E181").Value = 6

# Row 183: Aruba
$ws.Range("B183").Value = 118
$ws.Range("C183").Value = 1
$ws.Range("E183").Value = 14

# Row 210: Islas Malvinas
$ws.Range("A210").Value = 'Islas Malvinas'

# Row 211: Groenlandia
$ws.Range("A211").Value = 'Groenlandia'
